$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-08 05:48:15"
$ws.Range("N2").Value = "-4.3 °C 5:26 TU"
$ws.Range("O2").Value = "-2.4 °C"
$ws.Range("E3").Value = "2026-02-08 05:48:17"
$ws.Range("E4").Value = "2026-02-08 05:48:20"
$ws.Range("H4").Value = "70%"
$ws.Range("J4").Value = "1001.8 hPa"
$ws.Range("N4").Value = "6.1 °C 5:28 TU"
$ws.Range("O4").Value = "8.5 °C"
$ws.Range("E5").Value = "2026-02-08 05:48:22"
$ws.Range("N5").Value = "-5.1 °C 5:24 TU"
$ws.Range("E6").Value = "2026-02-08 05:48:25"
$ws.Range("H6").Value = "66%"
$ws.Range("J6").Value = "1001.6 hPa"
$ws.Range("N6").Value = "6.9 °C 5:24 TU"
$ws.Range("O6").Value = "8.5 °C"
$ws.Range("E7").Value = "2026-02-08 05:48:27"
$ws.Range("J7").Value = "1001.7 hPa"
$ws.Range("N7").Value = "10.3 °C 5:29 TU"
$ws.Range("O7").Value = "11.2 °C"
$ws.Range("E8").Value = "2026-02-08 05:48:29"
$ws.Range("H8").Value = "90%"
$ws.Range("J8").Value = "1001.7 hPa"
$ws.Range("N8").Value = "6.9 °C 5:29 TU"
$ws.Range("O8").Value = "8.2 °C"
$ws.Range("E9").Value = "2026-02-08 05:48:32"
$ws.Range("O9").Value = "7.5 °C"
$ws.Range("E10").Value = "2026-02-08 05:48:34"
$ws.Range("H10").Value = "93%"
$ws.Range("N10").Value = "5.0 °C 5:18 TU"
$ws.Range("O10").Value = "7.3 °C"
$ws.Range("E11").Value = "2026-02-08 05:48:36"
$ws.Range("N11").Value = "0.6 °C 5:00 TU"
$ws.Range("O11").Value = "1.6 °C"
$ws.Range("E12").Value = "2026-02-08 05:48:39"
$ws.Range("E13").Value = "2026-02-08 05:48:41"
$ws.Range("J13").Value = "1004.2 hPa"
$ws.Range("N13").Value = "-0.7 °C 5:25 TU"
$ws.Range("O13").Value = "0.7 °C"
$ws.Range("E14").Value = "2026-02-08 05:48:43"
$ws.Range("H14").Value = "95%"
$ws.Range("E15").Value = "2026-02-08 05:48:46"
$ws.Range("H15").Value = "85%"
$ws.Range("O15").Value = "6.1 °C"
$ws.Range("E16").Value = "2026-02-08 05:48:48"
$ws.Range("G16").Value = "70 cm"
$ws.Range("H16").Value = "86%"
$ws.Range("L16").Value = "24.8 km/h - 204º 5:08 TU"
$ws.Range("O16").Value = "-5.9 °C"
$ws.Range("E17").Value = "2026-02-08 05:48:50"
$ws.Range("E18").Value = "2026-02-08 05:48:53"
$ws.Range("H18").Value = "76%"
$ws.Range("J18").Value = "1002.0 hPa"
$ws.Range("N18").Value = "6.4 °C 5:29 TU"
$ws.Range("O18").Value = "8.0 °C"
$ws.Range("E19").Value = "2026-02-08 05:48:55"
$ws.Range("E20").Value = "2026-02-08 05:48:58"
$ws.Range("N20").Value = "-5.5 °C 5:29 TU"
$ws.Range("E21").Value = "2026-02-08 05:49:00"
$ws.Range("J21").Value = "1003.8 hPa"
$ws.Range("E22").Value = "2026-02-08 05:49:02"
$ws.Range("H22").Value = "94%"
$ws.Range("N22").Value = "-7.7 °C 5:29 TU"
$ws.Range("O22").Value = "-6.5 °C"
$ws.Range("E23").Value = "2026-02-08 05:49:05"
$ws.Range("H23").Value = "91%"
$ws.Range("E24").Value = "2026-02-08 05:49:07"
$ws.Range("H24").Value = "88%"
$ws.Range("N24").Value = "4.1 °C 5:25 TU"
$ws.Range("O24").Value = "6.8 °C"
$ws.Range("E25").Value = "2026-02-08 05:49:10"
$ws.Range("H25").Value = "87%"
$ws.Range("E26").Value = "2026-02-08 05:49:12"
$ws.Range("H26").Value = "76%"
$ws.Range("J26").Value = "1001.3 hPa"
$ws.Range("O26").Value = "1.5 °C"
$ws.Range("E27").Value = "2026-02-08 05:49:14"
$ws.Range("N27").Value = "-4.7 °C 5:13 TU"
$ws.Range("O27").Value = "-4.1 °C"
$ws.Range("E28").Value = "2026-02-08 05:49:17"
$ws.Range("H28").Value = "81%"
$ws.Range("J28").Value = "1001.9 hPa"
$ws.Range("N28").Value = "4.5 °C 5:15 TU"
$ws.Range("O28").Value = "5.7 °C"
$ws.Range("E29").Value = "2026-02-08 05:49:19"
$ws.Range("N29").Value = "7.9 °C 5:29 TU"
$ws.Range("O29").Value = "9.3 °C"
$ws.Range("E30").Value = "2026-02-08 05:49:21"
$ws.Range("H30").Value = "62%"
$ws.Range("J30").Value = "1001.1 hPa"
$ws.Range("N30").Value = "7.4 °C 5:29 TU"
$ws.Range("O30").Value = "9.5 °C"
$ws.Range("E31").Value = "2026-02-08 05:49:24"
$ws.Range("H31").Value = "58%"
$ws.Range("N31").Value = "9.4 °C 5:27 TU"
$ws.Range("O31").Value = "10.1 °C"
$ws.Range("E32").Value = "2026-02-08 05:49:27"
$ws.Range("O32").Value = "1.8 °C"
$ws.Range("E33").Value = "2026-02-08 05:49:29"
$ws.Range("J33").Value = "1003.5 hPa"
$ws.Range("N33").Value = "-0.5 °C 5:22 TU"
$ws.Range("O33").Value = "0.6 °C"
$ws.Range("E34").Value = "2026-02-08 05:49:31"
$ws.Range("H34").Value = "76%"
$ws.Range("O34").Value = "-1.1 °C"
$ws.Range("E35").Value = "2026-02-08 05:49:34"
$ws.Range("K35").Value = "-0.1 MJ/m2"
$ws.Range("N35").Value = "1.7 °C 5:25 TU"
$ws.Range("O35").Value = "3.9 °C"
$ws.Range("E36").Value = "2026-02-08 05:49:36"
$ws.Range("H36").Value = "67%"
$ws.Range("J36").Value = "1001.7 hPa"
$ws.Range("N36").Value = "9.9 °C 5:00 TU"
$ws.Range("O36").Value = "10.8 °C"
$ws.Range("E37").Value = "2026-02-08 05:49:38"
$ws.Range("J37").Value = "1003.1 hPa"
$ws.Range("E38").Value = "2026-02-08 05:49:41"
$ws.Range("H38").Value = "82%"
$ws.Range("N38").Value = "5.4 °C 5:26 TU"
$ws.Range("O38").Value = "7.6 °C"
$ws.Range("E39").Value = "2026-02-08 05:49:43"
$ws.Range("E40").Value = "2026-02-08 05:49:45"
$ws.Range("J40").Value = "1004.2 hPa"
$ws.Range("O40").Value = "2.7 °C"
$ws.Range("E41").Value = "2026-02-08 05:49:48"
$ws.Range("J41").Value = "1001.4 hPa"
$ws.Range("N41").Value = "7.1 °C 5:19 TU"
$ws.Range("O41").Value = "9.0 °C"
$ws.Range("E42").Value = "2026-02-08 05:49:50"
$ws.Range("H42").Value = "88%"
$ws.Range("N42").Value = "7.3 °C 5:29 TU"
$ws.Range("O42").Value = "9.5 °C"
$ws.Range("E43").Value = "2026-02-08 05:49:52"
$ws.Range("N43").Value = "2.6 °C 5:20 TU"
$ws.Range("O43").Value = "4.6 °C"
$ws.Range("E44").Value = "2026-02-08 05:49:54"
$ws.Range("N44").Value = "-6.2 °C 5:26 TU"
$ws.Range("O44").Value = "-5.6 °C"
$ws.Range("E45").Value = "2026-02-08 05:49:57"
$ws.Range("H45").Value = "66%"
$ws.Range("N45").Value = "0.1 °C 5:26 TU"
$ws.Range("O45").Value = "2.9 °C"
$ws.Range("E46").Value = "2026-02-08 05:49:59"
$ws.Range("H46").Value = "85%"
$ws.Range("J46").Value = "1002.0 hPa"
$ws.Range("N46").Value = "3.2 °C 5:06 TU"
$ws.Range("O46").Value = "6.7 °C"
